# Auto-generated edit script for horarios workbook update
# Applies the 2026-01-21 08:39:38 scrape refresh across all three sheets

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 08:39:38'
$ws.Range('A3').Value = 'Total filas: 110'
$ws.Range('A24').Value = '03:42:43'
$ws.Range('B24').Value = '05:35'
$ws.Range('C24').Value = '14_ABASTO'
$ws.Range('D24').Value = 113
$ws.Range('E24').Value = 'LP1912'
$ws.Range('A25').Value = '04:17:03'
$ws.Range('B25').Value = '05:35'
$ws.Range('C25').Value = '215B_EL PATO'
$ws.Range('D25').Value = 78
$ws.Range('E25').Value = 'LP1912'
$ws.Range('A55').Value = '07:17:57'
$ws.Range('B55').Value = '07:31'
$ws.Range('C55').Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range('D55').Value = 14
$ws.Range('E55').Value = 'LP1912'
$ws.Range('A56').Value = '07:17:57'
$ws.Range('B56').Value = '07:31'
$ws.Range('C56').Value = '16_SANTA ANA'
$ws.Range('D56').Value = 14
$ws.Range('E56').Value = 'LP1912'
$ws.Range('A57').Value = '07:17:57'
$ws.Range('B57').Value = '07:31'
$ws.Range('C57').Value = '11_ETCHEVERRY'
$ws.Range('D57').Value = 14
$ws.Range('E57').Value = 'LP1912'
$ws.Range('A83').Value = '08:39:38'
$ws.Range('B83').Value = '08:42'
$ws.Range('C83').Value = '81_EL PELIGRO'
$ws.Range('D83').Value = 3
$ws.Range('E83').Value = 'LP1912'
$ws.Range('A84').Value = '08:39:38'
$ws.Range('B84').Value = '08:43'
$ws.Range('C84').Value = '14_ABASTO'
$ws.Range('D84').Value = 4
$ws.Range('E84').Value = 'LP1912'
$ws.Range('A85').Value = '07:17:57'
$ws.Range('B85').Value = '08:53'
$ws.Range('C85').Value = '17_ROMERO'
$ws.Range('D85').Value = 96
$ws.Range('E85').Value = 'LP1912'
$ws.Range('A86').Value = '07:50:23'
$ws.Range('B86').Value = '08:53'
$ws.Range('C86').Value = '10_OLMOS'
$ws.Range('D86').Value = 63
$ws.Range('E86').Value = 'LP1912'
$ws.Range('A87').Value = '08:39:38'
$ws.Range('B87').Value = '08:54'
$ws.Range('C87').Value = '17_ROMERO'
$ws.Range('D87').Value = 15
$ws.Range('E87').Value = 'LP1912'
$ws.Range('A88').Value = '08:39:38'
$ws.Range('B88').Value = '08:55'
$ws.Range('C88').Value = '10_OLMOS'
$ws.Range('D88').Value = 16
$ws.Range('E88').Value = 'LP1912'
$ws.Range('A89').Value = '08:39:38'
$ws.Range('B89').Value = '09:01'
$ws.Range('C89').Value = '215A_EL PATO'
$ws.Range('D89').Value = 22
$ws.Range('E89').Value = 'LP1912'
$ws.Range('A90').Value = '08:39:38'
$ws.Range('B90').Value = '09:03'
$ws.Range('C90').Value = '11_ETCHEVERRY'
$ws.Range('D90').Value = 24
$ws.Range('E90').Value = 'LP1912'
$ws.Range('A91').Value = '08:39:38'
$ws.Range('B91').Value = '09:04'
$ws.Range('C91').Value = '23_HERNANDEZ'
$ws.Range('D91').Value = 25
$ws.Range('E91').Value = 'LP1912'
$ws.Range('A92').Value = '08:16:28'
$ws.Range('B92').Value = '09:08'
$ws.Range('C92').Value = '23_HERNANDEZ'
$ws.Range('D92').Value = 52
$ws.Range('E92').Value = 'LP1912'
$ws.Range('A93').Value = '08:39:38'
$ws.Range('B93').Value = '09:10'
$ws.Range('C93').Value = '16_P MOR-SANTA ANA'
$ws.Range('D93').Value = 31
$ws.Range('E93').Value = 'LP1912'
$ws.Range('A94').Value = '08:16:28'
$ws.Range('B94').Value = '09:13'
$ws.Range('C94').Value = '10_OLMOS'
$ws.Range('D94').Value = 57
$ws.Range('E94').Value = 'LP1912'
$ws.Range('A95').Value = '08:39:38'
$ws.Range('B95').Value = '09:16'
$ws.Range('C95').Value = '27_EL RETIRO'
$ws.Range('D95').Value = 37
$ws.Range('E95').Value = 'LP1912'
$ws.Range('A96').Value = '07:50:23'
$ws.Range('B96').Value = '09:17'
$ws.Range('C96').Value = '27_EL RETIRO'
$ws.Range('D96').Value = 87
$ws.Range('E96').Value = 'LP1912'
$ws.Range('A97').Value = '08:39:38'
$ws.Range('B97').Value = '09:21'
$ws.Range('C97').Value = '26_HERNANDEZ'
$ws.Range('D97').Value = 42
$ws.Range('E97').Value = 'LP1912'
$ws.Range('A98').Value = '08:39:38'
$ws.Range('B98').Value = '09:22'
$ws.Range('C98').Value = '16_SANTA ANA'
$ws.Range('D98').Value = 43
$ws.Range('E98').Value = 'LP1912'
$ws.Range('A99').Value = '08:39:38'
$ws.Range('B99').Value = '09:22'
$ws.Range('C99').Value = '17_ROMERO'
$ws.Range('D99').Value = 43
$ws.Range('E99').Value = 'LP1912'
$ws.Range('A100').Value = '08:39:38'
$ws.Range('B100').Value = '09:23'
$ws.Range('C100').Value = '11_ETCHEVERRY'
$ws.Range('D100').Value = 44
$ws.Range('E100').Value = 'LP1912'
$ws.Range('A101').Value = '07:50:23'
$ws.Range('B101').Value = '09:23'
$ws.Range('C101').Value = '17_ROMERO'
$ws.Range('D101').Value = 93
$ws.Range('E101').Value = 'LP1912'
$ws.Range('A102').Value = '08:16:28'
$ws.Range('B102').Value = '09:29'
$ws.Range('C102').Value = '16_SANTA ANA'
$ws.Range('D102').Value = 73
$ws.Range('E102').Value = 'LP1912'
$ws.Range('A103').Value = '07:50:23'
$ws.Range('B103').Value = '09:31'
$ws.Range('C103').Value = '16_SANTA ANA'
$ws.Range('D103').Value = 101
$ws.Range('E103').Value = 'LP1912'
$ws.Range('A104').Value = '08:39:38'
$ws.Range('B104').Value = '09:32'
$ws.Range('C104').Value = '15_ABASTO'
$ws.Range('D104').Value = 53
$ws.Range('E104').Value = 'LP1912'
$ws.Range('A105').Value = '08:39:38'
$ws.Range('B105').Value = '09:33'
$ws.Range('C105').Value = '10_OLMOS'
$ws.Range('D105').Value = 54
$ws.Range('E105').Value = 'LP1912'
$ws.Range('A106').Value = '08:39:38'
$ws.Range('B106').Value = '09:34'
$ws.Range('C106').Value = '16_SANTA ANA'
$ws.Range('D106').Value = 55
$ws.Range('E106').Value = 'LP1912'
$ws.Range('A107').Value = '08:39:38'
$ws.Range('B107').Value = '09:39'
$ws.Range('C107').Value = '23_HERNANDEZ'
$ws.Range('D107').Value = 60
$ws.Range('E107').Value = 'LP1912'
$ws.Range('A108').Value = '08:39:38'
$ws.Range('B108').Value = '09:41'
$ws.Range('C108').Value = '215C_EL PATO'
$ws.Range('D108').Value = 62
$ws.Range('E108').Value = 'LP1912'
$ws.Range('A109').Value = '08:39:38'
$ws.Range('B109').Value = '09:42'
$ws.Range('C109').Value = '10_OLMOS'
$ws.Range('D109').Value = 63
$ws.Range('E109').Value = 'LP1912'
$ws.Range('A110').Value = '08:16:28'
$ws.Range('B110').Value = '09:42'
$ws.Range('C110').Value = '215C_EL PATO'
$ws.Range('D110').Value = 86
$ws.Range('E110').Value = 'LP1912'
$ws.Range('A111').Value = '08:39:38'
$ws.Range('B111').Value = '09:43'
$ws.Range('C111').Value = '14_ABASTO'
$ws.Range('D111').Value = 64
$ws.Range('E111').Value = 'LP1912'
$ws.Range('A112').Value = '08:39:38'
$ws.Range('B112').Value = '10:10'
$ws.Range('C112').Value = '16_P MOR-SANTA ANA'
$ws.Range('D112').Value = 91
$ws.Range('E112').Value = 'LP1912'
$ws.Range('A113').Value = '08:39:38'
$ws.Range('B113').Value = '10:12'
$ws.Range('C113').Value = '15_ABASTO'
$ws.Range('D113').Value = 93
$ws.Range('E113').Value = 'LP1912'
$ws.Range('A114').Value = '08:39:38'
$ws.Range('B114').Value = '10:21'
$ws.Range('C114').Value = '26_HERNANDEZ'
$ws.Range('D114').Value = 102
$ws.Range('E114').Value = 'LP1912'
$ws.Range('A115').Value = '08:39:38'
$ws.Range('B115').Value = '10:26'
$ws.Range('C115').Value = '215A_EL PATO'
$ws.Range('D115').Value = 107
$ws.Range('E115').Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 08:39:38'
$ws.Range('A3').Value = 'Total filas: 22'
$ws.Range('A24').Value = '08:39:38'
$ws.Range('B24').Value = '09:01'
$ws.Range('C24').Value = '215A_EL PATO'
$ws.Range('D24').Value = 22
$ws.Range('E24').Value = 'LP1912'
$ws.Range('A25').Value = '08:39:38'
$ws.Range('B25').Value = '09:41'
$ws.Range('C25').Value = '215C_EL PATO'
$ws.Range('D25').Value = 62
$ws.Range('E25').Value = 'LP1912'
$ws.Range('A26').Value = '08:16:28'
$ws.Range('B26').Value = '09:42'
$ws.Range('C26').Value = '215C_EL PATO'
$ws.Range('D26').Value = 86
$ws.Range('E26').Value = 'LP1912'
$ws.Range('A27').Value = '08:39:38'
$ws.Range('B27').Value = '10:26'
$ws.Range('C27').Value = '215A_EL PATO'
$ws.Range('D27').Value = 107
$ws.Range('E27').Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 08:39:38'
$ws.Range('A3').Value = 'Total filas: 25'
$ws.Range('A27').Value = '08:39:38'
$ws.Range('B27').Value = '09:08'
$ws.Range('C27').Value = '215D_LA PLATA'
$ws.Range('D27').Value = 29
$ws.Range('E27').Value = 'L6203'
$ws.Range('A29').Value = '08:39:38'
$ws.Range('B29').Value = '10:02'
$ws.Range('C29').Value = '215B_LP-P MOR-40 Y 115'
$ws.Range('D29').Value = 83
$ws.Range('E29').Value = 'L6173'
$ws.Range('A30').Value = '08:16:28'
$ws.Range('B30').Value = '10:03'
$ws.Range('C30').Value = '215B_LP-P MOR-40 Y 115'
$ws.Range('D30').Value = 107
$ws.Range('E30').Value = 'L6173'
